# Auto-generated cell updates applying the Sheets diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1811.7142  # was 1513.8077
$ws.Range("I33").Value = 1759.6875  # was 1403.238
$ws.Range("K33").Value = 1759.6875  # was 1403.238
$ws.Range("M33").Value = -1530.6875  # was -1174.238
$ws.Range("H70").Value = 2271.92  # was 2287.76
$ws.Range("I70").Value = 2065.6667  # was 2109.6667
$ws.Range("K70").Value = 6197.000100000001  # was 6329.000100000001
$ws.Range("M70").Value = -5927.000100000001  # was -6059.000100000001
$ws.Range("H73").Value = 2271.92  # was 2287.76
$ws.Range("I73").Value = 2065.6667  # was 2109.6667
$ws.Range("K73").Value = 6197.000100000001  # was 6329.000100000001
$ws.Range("M73").Value = -5261.000100000001  # was -5393.000100000001
$ws.Range("H86").Value = 3357172.8  # was 3355006
$ws.Range("I86").Value = 4024607.5  # was 3355006
$ws.Range("J86").Value = 20000  # was 0
$ws.Range("K86").Value = 4024607.5  # was 3355006
$ws.Range("L86").Value = 20000  # was 0
$ws.Range("M86").Value = -4023484.5  # was -3353883
$ws.Range("N86").Value = -22246  # newly added cell
$ws.Range("H89").Value = 3357172.8  # was 3355006
$ws.Range("I89").Value = 4024607.5  # was 3355006
$ws.Range("J89").Value = 20000  # was 0
$ws.Range("K89").Value = 20123037.5  # was 16775030
$ws.Range("L89").Value = 100000  # was 0
$ws.Range("M89").Value = -20117421.5  # was -16769414
$ws.Range("N89").Value = -111232  # newly added cell
$ws.Range("H100").Value = 1729  # was 1677.4783
$ws.Range("I100").Value = 1061.1765  # was 1032.4445
$ws.Range("K100").Value = 1061.1765  # was 1032.4445
$ws.Range("M100").Value = -520.1765  # was -491.4445000000001
$ws.Range("H107").Value = 246.80952  # was 244.18182
$ws.Range("I107").Value = 238.55556  # was 235.94737
$ws.Range("K107").Value = 238.55556  # was 235.94737
$ws.Range("M107").Value = 1681.44444  # was 1684.05263
$ws.Range("H132").Value = 1415.7255  # was 1856.2368
$ws.Range("I132").Value = 1440.04  # was 1948.5278
$ws.Range("J132").Value = 200  # was 195
$ws.Range("K132").Value = 4320.12  # was 5845.5834
$ws.Range("L132").Value = 600  # was 585
$ws.Range("M132").Value = -1790.12  # was -3315.5834
$ws.Range("N132").Value = -5660  # was -5645
$ws.Range("H137").Value = 2975.9062  # was 3007.0908
$ws.Range("I137").Value = 1984.9524  # was 2138.762
$ws.Range("J137").Value = 4867.727  # was 4526.6665
$ws.Range("K137").Value = 5954.857199999999  # was 6416.286
$ws.Range("L137").Value = 14603.181  # was 13579.9995
$ws.Range("M137").Value = -3404.857199999999  # was -3866.286
$ws.Range("N137").Value = -19703.181  # was -18679.9995
$ws.Range("H138").Value = 2628.5305  # was 2522.8777
$ws.Range("J138").Value = 3346.4153  # was 3280.3508
$ws.Range("L138").Value = 10039.2459  # was 9841.0524
$ws.Range("N138").Value = -20319.2459  # was -20121.0524

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9094170  # was 9807338
$ws.Range("I32").Value = 10206539  # was 10640858
$ws.Range("J32").Value = 9816  # was 13474.5
$ws.Range("K32").Value = 10206539  # was 10640858
$ws.Range("L32").Value = 9816  # was 13474.5
$ws.Range("M32").Value = -10206252  # was -10640571
$ws.Range("N32").Value = -10390  # was -14048.5
$ws.Range("H61").Value = 10898180  # was 10444215
$ws.Range("I61").Value = 15155111  # was 14709520
$ws.Range("J61").Value = 92124.84  # was 85615.92999999999
$ws.Range("K61").Value = 15155111  # was 14709520
$ws.Range("L61").Value = 92124.84  # was 85615.92999999999
$ws.Range("M61").Value = -15154899  # was -14709308
$ws.Range("N61").Value = -92548.84  # was -86039.92999999999
$ws.Range("H105").Value = 50185  # was 42456.668
$ws.Range("J105").Value = 70370  # was 48685
$ws.Range("L105").Value = 70370  # was 48685
$ws.Range("N105").Value = -77358  # was -55673
$ws.Range("H107").Value = 108500  # was 108000
$ws.Range("J107").Value = 108500  # was 108000
$ws.Range("L107").Value = 108500  # was 108000
$ws.Range("N107").Value = -116180  # was -115680
$ws.Range("H110").Value = 1084.5  # was 1124.5238
$ws.Range("I110").Value = 1104.05  # was 1149.3158
$ws.Range("K110").Value = 1104.05  # was 1149.3158
$ws.Range("M110").Value = 940.95  # was 895.6841999999999
$ws.Range("H136").Value = 10898180  # was 10444215
$ws.Range("I136").Value = 15155111  # was 14709520
$ws.Range("J136").Value = 92124.84  # was 85615.92999999999
$ws.Range("K136").Value = 45465333  # was 44128560
$ws.Range("L136").Value = 276374.52  # was 256847.79
$ws.Range("M136").Value = -45462783  # was -44126010
$ws.Range("N136").Value = -281474.52  # was -261947.79

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3181.4  # was 3178.5
$ws.Range("I86").Value = 3021.25  # was 3116.6365
$ws.Range("J86").Value = 3822  # was 3405.3333
$ws.Range("K86").Value = 3021.25  # was 3116.6365
$ws.Range("L86").Value = 3822  # was 3405.3333
$ws.Range("M86").Value = -1898.25  # was -1993.6365
$ws.Range("N86").Value = -6068  # was -5651.3333
$ws.Range("H89").Value = 3181.4  # was 3178.5
$ws.Range("I89").Value = 3021.25  # was 3116.6365
$ws.Range("J89").Value = 3822  # was 3405.3333
$ws.Range("K89").Value = 15106.25  # was 15583.1825
$ws.Range("L89").Value = 19110  # was 17026.6665
$ws.Range("M89").Value = -9490.25  # was -9967.182500000001
$ws.Range("N89").Value = -30342  # was -28258.6665
$ws.Range("H107").Value = 1601.8064  # was 1648.5333
$ws.Range("I107").Value = 1146.1786  # was 1181.2222
$ws.Range("K107").Value = 1146.1786  # was 1181.2222
$ws.Range("M107").Value = 773.8214  # was 738.7778000000001
$ws.Range("H126").Value = 34999  # was 34499.5
$ws.Range("J126").Value = 34999  # was 34499.5
$ws.Range("L126").Value = 34999  # was 34499.5
$ws.Range("N126").Value = -44879  # was -44379.5
$ws.Range("H134").Value = 22935.064  # was 24037.092
$ws.Range("I134").Value = 1012.2195  # was 1105.35
$ws.Range("J134").Value = 202702.4  # was 253354.5
$ws.Range("K134").Value = 3036.6585  # was 3316.05
$ws.Range("L134").Value = 608107.2  # was 760063.5
$ws.Range("M134").Value = -501.6585  # was -781.0499999999997
$ws.Range("N134").Value = -613177.2  # was -765133.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 788532.9  # was 788585.4399999999
$ws.Range("I31").Value = 17249  # was 17392.363
$ws.Range("K31").Value = 17249  # was 17392.363
$ws.Range("M31").Value = -16954  # was -17097.363
$ws.Range("H34").Value = 788532.9  # was 788585.4399999999
$ws.Range("I34").Value = 17249  # was 17392.363
$ws.Range("K34").Value = 17249  # was 17392.363
$ws.Range("M34").Value = -17047  # was -17190.363
$ws.Range("H51").Value = 39666.668  # was 69000
$ws.Range("I51").Value = 25000  # was 0
$ws.Range("K51").Value = 25000  # was 0
$ws.Range("M51").Value = -24264  # newly added cell
$ws.Range("H60").Value = 28598  # was 26932
$ws.Range("J60").Value = 69997.5  # was 64999.5
$ws.Range("L60").Value = 69997.5  # was 64999.5
$ws.Range("N60").Value = -71019.5  # was -66021.5
$ws.Range("H61").Value = 39666.668  # was 69000
$ws.Range("I61").Value = 25000  # was 0
$ws.Range("K61").Value = 25000  # was 0
$ws.Range("M61").Value = -24652  # newly added cell
$ws.Range("H99").Value = 3421.4285  # was 3735.6428
$ws.Range("I99").Value = 3418.182  # was 3724.9167
$ws.Range("J99").Value = 3433.3333  # was 3800
$ws.Range("K99").Value = 3418.182  # was 3724.9167
$ws.Range("L99").Value = 3433.3333  # was 3800
$ws.Range("M99").Value = -1920.182  # was -2226.9167
$ws.Range("N99").Value = -6429.3333  # was -6796
$ws.Range("H126").Value = 3421.4285  # was 3735.6428
$ws.Range("I126").Value = 3418.182  # was 3724.9167
$ws.Range("J126").Value = 3433.3333  # was 3800
$ws.Range("K126").Value = 10254.546  # was 11174.7501
$ws.Range("L126").Value = 10299.9999  # was 11400
$ws.Range("M126").Value = -7784.545999999998  # was -8704.750100000001
$ws.Range("N126").Value = -15239.9999  # was -16340
$ws.Range("H141").Value = 144496.67  # was 158797.8
$ws.Range("J141").Value = 144496.67  # was 158797.8
$ws.Range("L141").Value = 144496.67  # was 158797.8
$ws.Range("N141").Value = -154856.67  # was -169157.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I80").Value = 5500  # was 0
$ws.Range("K80").Value = 16500  # was 0
$ws.Range("M80").Value = -15564  # newly added cell
$ws.Range("I83").Value = 5500  # was 0
$ws.Range("K83").Value = 49500  # was 0
$ws.Range("M83").Value = -44820  # newly added cell
$ws.Range("H92").Value = 2004979.8  # was 1670831.5
$ws.Range("J92").Value = 4966.3335  # was 3747.25
$ws.Range("L92").Value = 14899.0005  # was 11241.75
$ws.Range("N92").Value = -17395.0005  # was -13737.75
$ws.Range("H129").Value = 30391392  # was 27858850
$ws.Range("J129").Value = 55716384  # was 47757028
$ws.Range("L129").Value = 167149152  # was 143271084
$ws.Range("N129").Value = -167159152  # was -143281084

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 64140.6  # was 89989.5
$ws.Range("I5").Value = 64140.6  # was 106651
$ws.Range("J5").Value = 0  # was 40005
$ws.Range("K5").Value = 64140.6  # was 106651
$ws.Range("L5").Value = 0  # was 40005
$ws.Range("M5").Value = -64028.6  # was -106539
$ws.Range("N5").ClearContents()  # was -40229
$ws.Range("H109").Value = 62642.5  # was 0
$ws.Range("J109").Value = 62642.5  # was 0
$ws.Range("L109").Value = 62642.5  # was 0
$ws.Range("N109").Value = -64722.5  # newly added cell
$ws.Range("H122").Value = 2302  # was 2461
$ws.Range("I122").Value = 1935.3334  # was 2149.5
$ws.Range("K122").Value = 5806.0002  # was 6448.5
$ws.Range("M122").Value = -3356.0002  # was -3998.5
$ws.Range("H126").Value = 3346.5557  # was 3732
$ws.Range("I126").Value = 3412.2  # was 4197.25
$ws.Range("J126").Value = 3264.5  # was 3266.75
$ws.Range("K126").Value = 10236.6  # was 12591.75
$ws.Range("L126").Value = 9793.5  # was 9800.25
$ws.Range("M126").Value = -7766.599999999999  # was -10121.75
$ws.Range("N126").Value = -14733.5  # was -14740.25
$ws.Range("H133").Value = 75000  # was 72500
$ws.Range("J133").Value = 50000  # was 63333.332
$ws.Range("L133").Value = 50000  # was 63333.332
$ws.Range("N133").Value = -60120  # was -73453.33199999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 29373.5  # was 26666.334
$ws.Range("J48").Value = 37495  # was 0
$ws.Range("L48").Value = 37495  # was 0
$ws.Range("N48").Value = -38817  # newly added cell

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 68339.37  # was 73223.3
$ws.Range("J4").Value = 37716.375  # was 40318.715
$ws.Range("L4").Value = 37716.375  # was 40318.715
$ws.Range("N4").Value = -37942.375  # was -40544.715
$ws.Range("H51").Value = 37999  # was 39332.668
$ws.Range("J51").Value = 37999  # was 39332.668
$ws.Range("L51").Value = 37999  # was 39332.668
$ws.Range("N51").Value = -39019  # was -40352.668
$ws.Range("H107").Value = 62501176  # was 55556656
$ws.Range("J107").Value = 555  # was 527.5
$ws.Range("L107").Value = 1665  # was 1582.5
$ws.Range("N107").Value = -5505  # was -5422.5
$ws.Range("H119").Value = 91029.336  # was 91062.664
$ws.Range("J119").Value = 91029.336  # was 91062.664
$ws.Range("L119").Value = 91029.336  # was 91062.664
$ws.Range("N119").Value = -100705.336  # was -100738.664
$ws.Range("H138").Value = 58333.332  # was 60000
$ws.Range("J138").Value = 58333.332  # was 60000
$ws.Range("L138").Value = 58333.332  # was 60000
$ws.Range("N138").Value = -68613.33199999999  # was -70280
